$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The post "「重力波の発見は新しい窓を開いた」" (row 543) was removed.
# Deleting the entire row shifts all subsequent rows up by one, matching
# the rest of the diff (rows 544..603 -> 543..602).
$ws.Rows("543").Delete()
